$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, bordered, centered) from an existing header cell (AC1)
# onto the three new header cells so they match the rest of the header row.
$headerStyleSrc = $ws.Range("AC1")
$newHeaders = $ws.Range("AD1:AF1")
$headerStyleSrc.Copy()
$newHeaders.PasteSpecial(-4122)

# Set the new header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every player row.
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 82
    $ws.Cells.Item($r, 31).Value = 80
    $ws.Cells.Item($r, 32).Value = 0
}
